# Automatische test-sync: 2025-08-06 20:35:50
# Appends the newest mail-log entry to the "Logs" sheet (row 18) and bumps
# the "Retour / Terugbetaling" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 18

$logs.Range("A$newRow").Value = "Weten jullie al iets over mijn retour?"
$logs.Range("B$newRow").Value = "mailmind.test@zohomail.eu"
$logs.Range("C$newRow").Value = "Testmail #1: Weten jullie al iets over mijn retour?"
$logs.Range("D$newRow").Value = "Retour / Terugbetaling"
$logs.Range("E$newRow").Value = "Beste klant,`r`nBedankt voor je e-mail. Om je vraag over de status van je retour te kunnen beantwoorden, heb ik wat meer informatie nodig. Zou je alsjeblieft je ordernummer of traceernummer van de retourzending kunnen doorgeven? Dan kan ik dit voor je nakijken en je zo goed mogelijk helpen.`r`nMet vriendelijke groet,`r`n[Naam] E-mailassistent"
$logs.Range("F$newRow").Value = "2025-08-06 20:35:19"
$logs.Range("G$newRow").Value = "Ja"
$logs.Range("H$newRow").Value = "Nee"
$logs.Range("I$newRow").Value = "Ja"
$logs.Range("J$newRow").Value = "Nee"

# Re-fit the row height after the multi-line "Antwoord" text was set, so the
# row doesn't keep an explicit custom height like a real paste/entry wouldn't.
$logs.Rows.Item($newRow).AutoFit()

# Extend the conditional-formatting ranges so the new row inherits the
# same colour rules as the rest of the table (D/G/H/I/J columns).
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$($col)2:$($col)17")
    $newRange = $logs.Range("$($col)2:$($col)$newRow")
    $rules = $oldRange.FormatConditions
    for ($i = 1; $i -le $rules.Count; $i++) {
        $rules.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary: one more "Retour / Terugbetaling" mail.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B6").Value = 2
